$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # RUNMANAGER
$ws2 = $wb.Worksheets.Item(2)  # RETAIL_DATA

# --- RUNMANAGER: drop the "Priority" column (old column D), "Count" becomes new column D ---
$ws1.Columns.Item(4).Delete()

# Add new test case row 12 ("quickLinksCheck") to RUNMANAGER.
# Copy D11 -> D12 first so the new row's Count cell inherits the same text-number style (s="1")
$ws1.Range("D11").Copy($ws1.Range("D12"))
$ws1.Cells.Item(12, 1).Value = "quickLinksCheck"
$ws1.Cells.Item(12, 2).Value = "checking dashboard quick linnks"
$ws1.Cells.Item(12, 3).Value = "yes"

# --- RETAIL_DATA: add matching new test case row 13 ("quickLinksCheck") ---
# Copy full row 12 down to row 13 so format/styles (incl. hyperlink cell style) carry over
$ws2.Range("A12:E12").Copy($ws2.Range("A13:E13"))
$ws2.Cells.Item(13, 1).Value = "quickLinksCheck"
$ws2.Cells.Item(13, 4).Value = "spcbtest"

# Register the real hyperlink relationship for E13 (Copy above only cloned the visual style/text)
$style13 = $ws2.Cells.Item(13, 5).Style
$ws2.Hyperlinks.Add($ws2.Range("E13"), "mailto:Asdf@123")
$ws2.Cells.Item(13, 5).Style = $style13

# TC00 row in RETAIL_DATA uses "spcb" rather than "spcbtest"
$ws2.Cells.Item(3, 4).Value = "spcb"

# --- Active tab / selection bookkeeping: RETAIL_DATA becomes the active sheet ---
$ws1.Range("A1").Select()
$ws2.Activate()
$ws2.Range("B19").Select()
